$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new job posting row (Job_Id = 17) in row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = "Java Developer"
$ws.Cells.Item(18, 3).Value = "dfdsfsfdsfsd"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
